# The deck's theme (Integral) is replaced with the stock "Office Theme"
# colour scheme. (Font scheme / format scheme are already identical
# between the two themes in this deck, so only the 12 theme colours -
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - actually change.)

function HexToRgbVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToRgbVal $officeThemeColors[$i - 1]
}
